# Atualizado por script em 12-11-2023 14:45
#
# This script reproduces a data refresh on the CFL Group B (Czech Republic)
# betting-odds sheet:
#   - A handful of already-recorded matches had their row order corrected
#     (the match rows themselves did not move, only the F:V match-data
#     block was shuffled into the right row), and
#   - two brand-new matches were appended at the bottom (rows 119/120),
#     pushing the used range from A1:V118 to A1:V120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# --- Re-shuffle the match-data block (columns F:V) among existing rows. ---
# Columns A:E (Indice, pais, torneio, temporada, data_partida) stay put;
# only the home/away/odds/url block moves, matching the source diff.

# Pairwise swap: row 70 <-> row 71
$v70 = Get-RowValues 70
$v71 = Get-RowValues 71
Set-RowValues 70 $v71
Set-RowValues 71 $v70

# Pairwise swap: row 85 <-> row 86
$v85 = Get-RowValues 85
$v86 = Get-RowValues 86
Set-RowValues 85 $v86
Set-RowValues 86 $v85

# Pairwise swap: row 99 <-> row 100
$v99 = Get-RowValues 99
$v100 = Get-RowValues 100
Set-RowValues 99 $v100
Set-RowValues 100 $v99

# Pairwise swap: row 107 <-> row 108
$v107 = Get-RowValues 107
$v108 = Get-RowValues 108
Set-RowValues 107 $v108
Set-RowValues 108 $v107

# Three-way rotation: row 115 <- 116 <- 117 <- 115
$v115 = Get-RowValues 115
$v116 = Get-RowValues 116
$v117 = Get-RowValues 117
Set-RowValues 115 $v116
Set-RowValues 116 $v117
Set-RowValues 117 $v115

# --- Append two newly scraped matches as rows 119 and 120. ---
# Copy formatting (borders/number formats/etc.) down from the last
# existing data row so the new rows match the table's styling.
$ws.Range("A118:V118").Copy()
$ws.Range("A119:V120").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "czech-republic"
$ws.Range("C119").Value = "cfl-group-b"
$ws.Range("D119").Value = "2023-2024"
$ws.Range("E119").Value = 45242.42708333334
$ws.Range("F119").Value = "Pardubice B"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = "Zivanice"
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 2.08
$ws.Range("K119").Value = "12/11/2023 00:12"
$ws.Range("L119").Value = 1.62
$ws.Range("M119").Value = "12/11/2023 10:14"
$ws.Range("N119").Value = 3.5
$ws.Range("O119").Value = "12/11/2023 00:12"
$ws.Range("P119").Value = 4.28
$ws.Range("Q119").Value = "12/11/2023 10:14"
$ws.Range("R119").Value = 2.99
$ws.Range("S119").Value = "12/11/2023 00:12"
$ws.Range("T119").Value = 4.28
$ws.Range("U119").Value = "12/11/2023 10:14"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-b/pardubice-zivanice/b5mrhCg6/"

$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "czech-republic"
$ws.Range("C120").Value = "cfl-group-b"
$ws.Range("D120").Value = "2023-2024"
$ws.Range("E120").Value = 45242.45833333334
$ws.Range("F120").Value = "Teplice B"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Liberec B"
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = 2.5
$ws.Range("K120").Value = "12/11/2023 00:12"
$ws.Range("L120").Value = 2.41
$ws.Range("M120").Value = "12/11/2023 10:27"
$ws.Range("N120").Value = 3.35
$ws.Range("O120").Value = "12/11/2023 00:12"
$ws.Range("P120").Value = 3.71
$ws.Range("Q120").Value = "12/11/2023 10:27"
$ws.Range("R120").Value = 2.48
$ws.Range("S120").Value = "12/11/2023 00:12"
$ws.Range("T120").Value = 2.49
$ws.Range("U120").Value = "12/11/2023 10:27"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-b/teplice-liberec/tUhwghv0/"
